$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:D2").Value = 0
$ws.Range("B3:D3").Value = 0
$ws.Range("B4:C4").Value = 0
$ws.Range("D4").Value = 695.25
